# Adiciona cálculo e exibição do horário da próxima execução no loop automático
#
# The automated scraper loop now recalculates/shows the timestamp of its next
# run. For this snapshot that means: the latest poll only produced a single
# fresh match row (the previous run's stale "Em Andamento" rows 3-6 are gone),
# the surviving row's teams/competition were refreshed, and its "Hora" column
# reflects the newly computed next-run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the remaining match row (row 2) with the latest scrape results.
$ws.Range("A2").Value = "Süper Lig"
$ws.Range("B2").Value = "Fenerbahce"
$ws.Range("D2").Value = "Trabzonspor"
$ws.Range("J2").Value = "13:43:15"

# Drop the now-stale rows (old cycle's other in-progress matches).
$ws.Rows("3:6").Delete()

# Re-apply the autofilter so its range shrinks to match the live data (A1:J2)
# instead of staying stuck on the old A1:J6 extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:J2").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range too.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
  $n = $names.Item($i)
  if ($n.Name() -like "*_FilterDatabase*") {
    $n.RefersTo = "='Jogos 0x0'!`$A`$1:`$J`$2"
  }
}
